# Handles float input without breaking stuff
# Applies the marksheet update: scores in the summary block (rows 10-12),
# collapses the "Student Ans" duplicate block from three columns (A/B, D/E,
# G/H) down to two (A/B, D/E only for the first 3 rows), and promotes the
# previous "Correct Ans" value into the newly-freed "Student Ans" slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Summary block (rows 10-12): give the row-label cells in column A the
# same "mtitleStyle" formatting already used by the header row (row 9),
# and refresh the numeric / text results.
# ---------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 18
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 72
$ws.Range("E12").Value = "72/112"

# ---------------------------------------------------------------------
# Question answer grid (rows 15-40): drop the third "Student Ans /
# Correct Ans" pair (columns G:H) entirely, and drop the second pair
# (columns D:E) past row 18.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# The remaining D:E pair (rows 16-18) now mirrors the old "Correct Ans"
# text, formatted like the "Student Ans" column (style copied from B10,
# i.e. "correctStyle").
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# Column A ("Student Ans") picks up the previous "Correct Ans" (column B)
# value for the rows that actually hold a question; other rows are left
# untouched (blank, as before).
$aAnswers = [ordered]@{
  16 = "Option A"
  17 = "Option D"
  18 = "Option B"
  19 = "Option C"
  21 = "Option C"
  25 = "Option A"
  26 = "Option C"
  27 = "Option A"
  29 = "Option D"
  30 = "Option B"
  32 = "Option C"
  33 = "Option D"
  37 = "Option A"
  38 = "Option A"
  39 = "Option D"
}

$ws.Range("B10").Copy()
foreach ($r in $aAnswers.Keys) {
  $ws.Range("A$r").PasteSpecial(-4122)
}
foreach ($r in $aAnswers.Keys) {
  $ws.Range("A$r").Value = $aAnswers[$r]
}
